# Remove testing social media values
# The "values" sheet has placeholder "test" strings in column B for the
# optional variables (google_analytics_id, opengraph_image,
# opengraph_description, twitter_description) in rows 3-6. These test
# values should be cleared, leaving the cells blank (keeping formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("values")

$ws.Range("B3:B6").ClearContents()
